# Auto-generated edit script applying the Midgardsormr_Profits scheduled data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 5000
$ws.Range("I43").Value = 5000
$ws.Range("K43").Value = 5000
$ws.Range("M43").Value = -4931
$ws.Range("H111").Value = 2381.6875
$ws.Range("I111").Value = 1716
$ws.Range("J111").Value = 3237.5715
$ws.Range("K111").Value = 5148
$ws.Range("L111").Value = 9712.7145
$ws.Range("M111").Value = -2081
$ws.Range("N111").Value = -15846.7145
$ws.Range("H137").Value = 56545.273
$ws.Range("I137").Value = 100648.664
$ws.Range("K137").Value = 301945.992
$ws.Range("M137").Value = -299395.992
$ws.Range("H138").Value = 30139.445
$ws.Range("I138").Value = 1780.8182
$ws.Range("K138").Value = 5342.4546
$ws.Range("M138").Value = -202.4546

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 30078.475
$ws.Range("J32").Value = 8666
$ws.Range("L32").Value = 8666
$ws.Range("N32").Value = -9240
$ws.Range("H45").Value = 3577.5715
$ws.Range("I45").Value = 1449
$ws.Range("K45").Value = 1449
$ws.Range("M45").Value = -1072
$ws.Range("H61").Value = 23164
$ws.Range("I61").Value = 3606.6667
$ws.Range("K61").Value = 3606.6667
$ws.Range("M61").Value = -3394.6667
$ws.Range("H63").Value = 3520.3572
$ws.Range("I63").Value = 3160.875
$ws.Range("K63").Value = 3160.875
$ws.Range("M63").Value = -2474.875
$ws.Range("H66").Value = 3520.3572
$ws.Range("I66").Value = 3160.875
$ws.Range("K66").Value = 15804.375
$ws.Range("M66").Value = -12372.375
$ws.Range("H74").Value = 383442.3
$ws.Range("J74").Value = 15947.25
$ws.Range("L74").Value = 15947.25
$ws.Range("N74").Value = -17695.25
$ws.Range("H77").Value = 383442.3
$ws.Range("J77").Value = 15947.25
$ws.Range("L77").Value = 79736.25
$ws.Range("N77").Value = -88472.25
$ws.Range("H97").Value = 703.5217
$ws.Range("I97").Value = 617.76746
$ws.Range("K97").Value = 617.76746
$ws.Range("M97").Value = -121.76746
$ws.Range("H102").Value = 2753.1538
$ws.Range("I102").Value = 2664.7778
$ws.Range("J102").Value = 2952
$ws.Range("K102").Value = 2664.7778
$ws.Range("L102").Value = 2952
$ws.Range("M102").Value = -1042.7778
$ws.Range("N102").Value = -6196
$ws.Range("H136").Value = 23164
$ws.Range("I136").Value = 3606.6667
$ws.Range("K136").Value = 10820.0001
$ws.Range("M136").Value = -8270.000100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 26946.75
$ws.Range("H78").Value = 26946.75
$ws.Range("H86").Value = 2266.4443
$ws.Range("I86").Value = 2266.4443
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 2266.4443
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -1143.4443
$ws.Range("H89").Value = 2266.4443
$ws.Range("I89").Value = 2266.4443
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 11332.2215
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -5716.2215
$ws.Range("H105").Value = 1714.2593
$ws.Range("I105").Value = 1714.2593
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1714.2593
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 32.74070000000006
$ws.Range("H107").Value = 2161.1875
$ws.Range("I107").Value = 2161.1875
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2161.1875
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -241.1875
$ws.Range("N86").ClearContents()
$ws.Range("N89").ClearContents()
$ws.Range("N105").ClearContents()
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2070.125
$ws.Range("I16").Value = 2009.4615
$ws.Range("K16").Value = 2009.4615
$ws.Range("M16").Value = -1722.4615
$ws.Range("H31").Value = 8334252
$ws.Range("I31").Value = 8334252
$ws.Range("K31").Value = 8334252
$ws.Range("M31").Value = -8333957
$ws.Range("H34").Value = 8334252
$ws.Range("I34").Value = 8334252
$ws.Range("K34").Value = 8334252
$ws.Range("M34").Value = -8334050
$ws.Range("H58").Value = 2210.0588
$ws.Range("I58").Value = 1152.1818
$ws.Range("K58").Value = 1152.1818
$ws.Range("M58").Value = -949.1818000000001
$ws.Range("H86").Value = 73688.39999999999
$ws.Range("J86").Value = 17499.25
$ws.Range("L86").Value = 17499.25
$ws.Range("N86").Value = -19745.25
$ws.Range("H89").Value = 73688.39999999999
$ws.Range("J89").Value = 17499.25
$ws.Range("L89").Value = 87496.25
$ws.Range("N89").Value = -98728.25
$ws.Range("H99").Value = 6197.25
$ws.Range("I99").Value = 4852
$ws.Range("K99").Value = 4852
$ws.Range("M99").Value = -3354
$ws.Range("H113").Value = 2070.125
$ws.Range("I113").Value = 2009.4615
$ws.Range("K113").Value = 2009.4615
$ws.Range("M113").Value = 160.5385000000001
$ws.Range("H126").Value = 6197.25
$ws.Range("I126").Value = 4852
$ws.Range("K126").Value = 14556
$ws.Range("M126").Value = -12086
$ws.Range("H136").Value = 2210.0588
$ws.Range("I136").Value = 1152.1818
$ws.Range("K136").Value = 3456.5454
$ws.Range("M136").Value = -906.5454

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 4436.3706
$ws.Range("I59").Value = 3926.6667
$ws.Range("J59").Value = 4500.0835
$ws.Range("K59").Value = 11780.0001
$ws.Range("L59").Value = 13500.2505
$ws.Range("M59").Value = -11240.0001
$ws.Range("N59").Value = -14580.2505
$ws.Range("H124").Value = 8740.317999999999
$ws.Range("I124").Value = 2321.5
$ws.Range("J124").Value = 9382.200000000001
$ws.Range("K124").Value = 6964.5
$ws.Range("L124").Value = 28146.6
$ws.Range("M124").Value = -2054.5
$ws.Range("N124").Value = -37966.60000000001
$ws.Range("H137").Value = 2621.2666
$ws.Range("I137").Value = 2406.6155
$ws.Range("K137").Value = 7219.8465
$ws.Range("M137").Value = -2119.8465

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1975.2593
$ws.Range("I113").Value = 1974.2778
$ws.Range("J113").Value = 1977.2222
$ws.Range("K113").Value = 1974.2778
$ws.Range("L113").Value = 1977.2222
$ws.Range("M113").Value = 195.7221999999999
$ws.Range("N113").Value = -6317.2222

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 17490.092
$ws.Range("J45").Value = 17490.092
$ws.Range("L45").Value = 17490.092
$ws.Range("N45").Value = -18472.092
$ws.Range("H136").Value = 25930.24
$ws.Range("I136").Value = 29068.455
$ws.Range("J136").Value = 2916.6667
$ws.Range("K136").Value = 87205.36500000001
$ws.Range("L136").Value = 8750.000100000001
$ws.Range("M136").Value = -84655.36500000001
$ws.Range("N136").Value = -13850.0001
